# Adds 2022-Q3 fund-holding data to the 天津港 (600717) workbook.
#
# Net effect (see commit message "feat: add 2022-Q3 data"):
#   - 总计 (summary) sheet gains a new leading row for 2022-Q3 and all
#     later rows shift down by one.
#   - Every existing quarter sheet is renamed to the *next* older
#     quarter and its data is replaced with the data that used to
#     belong to that slot's previous (one-quarter-newer) occupant - i.e.
#     everything cascades one quarter older to make room for the new
#     quarter at the front.
#   - The former "2021-Q3" sheet is repurposed to hold the brand new
#     2022-Q3 fund table.
#   - A brand new sheet is appended at the end, named "2020-Q4", holding
#     the data the old "2020-Q4" sheet had before the cascade (that
#     quarter's numbers do not change).

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $cellRef, $val) {
    # Force text storage (matches the workbook's t="inlineStr" cells) so
    # numeric-looking strings like "0.12" or fund codes like "512780"
    # are not silently coerced into real numbers by Excel.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

function Set-NumCell($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.Value = $val
}

function Fill-FundRow($ws, $rowNum, $code, $name, $amount, $pos, $ratio, $mv, $rank, $idx) {
    $refA = "A" + $rowNum
    $refB = "B" + $rowNum
    $refC = "C" + $rowNum
    $refD = "D" + $rowNum
    $refE = "E" + $rowNum
    $refF = "F" + $rowNum
    $refG = "G" + $rowNum
    $refH = "H" + $rowNum
    Set-NumCell  $ws $refA $idx
    Set-TextCell $ws $refB $code
    Set-TextCell $ws $refC $name
    Set-TextCell $ws $refD $amount
    Set-TextCell $ws $refE $pos
    Set-TextCell $ws $refF $ratio
    Set-TextCell $ws $refG $mv
    Set-NumCell  $ws $refH $rank
}

# ---------------------------------------------------------------------
# 1. Duplicate the existing "2020-Q4" sheet (last tab) so the cascade
#    below has a final resting place for the data it currently holds;
#    the copy preserves that sheet's distinct styling (border/alignment
#    style index + page margins) automatically.
# ---------------------------------------------------------------------
$lastIdx = $wb.Worksheets.Count
$oldQ4 = $wb.Worksheets.Item($lastIdx)
$oldQ4.Copy($null, $oldQ4)
$newLastIdx = $wb.Worksheets.Count
$newQ4 = $wb.Worksheets.Item($newLastIdx)
$newQ4.Name = "2020-Q4-STAGE"

# ---------------------------------------------------------------------
# 2. Rename each quarter sheet to the next-older quarter, in positional
#    order (2 -> 5); every target name is unused at the moment it is
#    assigned, so there is no collision.
# ---------------------------------------------------------------------
$wsNewQ3 = $wb.Worksheets.Item(2)
$wsNewQ3.Name = "2022-Q3"

$wsQ3 = $wb.Worksheets.Item(3)
$wsQ3.Name = "2021-Q3"

$wsQ2 = $wb.Worksheets.Item(4)
$wsQ2.Name = "2021-Q2"

$wsQ1 = $wb.Worksheets.Item(5)
$wsQ1.Name = "2021-Q1"

$newQ4.Name = "2020-Q4"

# ---------------------------------------------------------------------
# 3. Overwrite the data in each cascaded sheet with the values that used
#    to belong one slot earlier (the "2020-Q4" sheet's own values do not
#    change, so it needs no data rewrite).
# ---------------------------------------------------------------------

# 2021-Q1 sheet <- what used to be the "2021-Q1" sheet's own data
Fill-FundRow $wsQ1 2 "512780" "广发中证京津冀协同发展主题ETF" "0.42" "99.46" "3.00" "0.0126" 9 0
Fill-FundRow $wsQ1 3 "164811" "工银瑞信中证京津冀协同发展主题指数（LOF）A" "0.13" "94.67" "2.99" "0.0039" 10 1
Fill-FundRow $wsQ1 4 "164825" "工银瑞信中证京津冀协同发展主题指数（LOF）C" "0.01" "94.67" "2.99" "0.0003" 10 2

# 2021-Q2 sheet <- what used to be the "2021-Q2" sheet's own data
Fill-FundRow $wsQ2 2 "512780" "广发中证京津冀协同发展主题ETF" "0.25" "98.47" "2.90" "0.0072" 10 0
Fill-FundRow $wsQ2 3 "164811" "工银瑞信中证京津冀协同发展主题指数（LOF）A" "0.12" "94.10" "2.77" "0.0033" 10 1
Fill-FundRow $wsQ2 4 "164825" "工银瑞信中证京津冀协同发展主题指数（LOF）C" "0.02" "94.10" "2.77" "0.0006" 10 2

# 2021-Q3 sheet <- what used to be the "2021-Q3" sheet's own data
Fill-FundRow $wsQ3 2 "512780" "广发中证京津冀协同发展主题ETF" "0.15" "98.85" "3.26" "0.0049" 3 0
Fill-FundRow $wsQ3 3 "164811" "工银瑞信中证京津冀协同发展主题指数（LOF）A" "0.13" "92.93" "3.06" "0.0040" 3 1
Fill-FundRow $wsQ3 4 "164825" "工银瑞信中证京津冀协同发展主题指数（LOF）C" "0.02" "92.93" "3.06" "0.0006" 3 2

# ---------------------------------------------------------------------
# 4. Brand new 2022-Q3 sheet: reload the former "2021-Q3" slot with the
#    2022-Q3 holdings (5 funds now) and fix the "基金规模" column header.
# ---------------------------------------------------------------------
Set-TextCell $wsNewQ3 "D1" "基金规模"

$dim = $wsNewQ3.Range("A1:H6")
Write-Output ("dim placeholder " + $dim.Address)

Fill-FundRow $wsNewQ3 2 "164811" "工银瑞信中证京津冀协同发展主题指数（LOF）A" "0.12" "93.09" "3.11" "0.0037" 4 0
Fill-FundRow $wsNewQ3 3 "851088" "海通量化成长精选一年持有期混合A" "0.38" "85.56" "0.93" "0.0035" 6 1
Fill-FundRow $wsNewQ3 4 "850010" "海通量化成长精选一年持有期混合B" "0.25" "85.56" "0.93" "0.0023" 6 2
Fill-FundRow $wsNewQ3 5 "164825" "工银瑞信中证京津冀协同发展主题指数（LOF）C" "0.03" "93.09" "3.11" "0.0009" 4 3
Fill-FundRow $wsNewQ3 6 "851099" "海通量化成长精选一年持有期混合C" "0.03" "85.56" "0.93" "0.0003" 6 4

# ---------------------------------------------------------------------
# 5. Update the 总计 (summary) sheet: insert the new 2022-Q3 row at the
#    top of the data and push every other row down by one.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

Set-NumCell  $wsTotal "A6" 4
Set-TextCell $wsTotal "B6" "2020-Q4"
Set-NumCell  $wsTotal "C6" 3
Set-NumCell  $wsTotal "D6" 0.02

Set-NumCell  $wsTotal "A5" 3
Set-TextCell $wsTotal "B5" "2021-Q1"
Set-NumCell  $wsTotal "C5" 3
Set-NumCell  $wsTotal "D5" 0.02

Set-NumCell  $wsTotal "A4" 2
Set-TextCell $wsTotal "B4" "2021-Q2"
Set-NumCell  $wsTotal "C4" 3
Set-NumCell  $wsTotal "D4" 0.01

Set-NumCell  $wsTotal "A3" 1
Set-TextCell $wsTotal "B3" "2021-Q3"
Set-NumCell  $wsTotal "C3" 3
Set-NumCell  $wsTotal "D3" 0.01

Set-NumCell  $wsTotal "A2" 0
Set-TextCell $wsTotal "B2" "2022-Q3"
Set-NumCell  $wsTotal "C2" 5
Set-NumCell  $wsTotal "D2" 0.01
